# Update the "dSF" (column F) values to match the repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -1
$ws.Range("F12").Value = 0
$ws.Range("F16").Value = -3
$ws.Range("F19").Value = -4
$ws.Range("F20").Value = 3
$ws.Range("F21").Value = 2
$ws.Range("F23").Value = 3
$ws.Range("F25").Value = 0
$ws.Range("F27").Value = -5
